$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.218.42"
$ws.Range("E2").Value = "  +3.98%  "
$ws.Range("D3").Value = "2.986.47"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.33%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "2.983.24"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "66.211.00"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "3.482.85"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "2.990.48"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +4.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.88%  "
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000103"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.89%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "395.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D47").Value = "2.739.66"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.07%  "
$ws.Range("E51").Value = "  +1.86%  "
